$d = $word.ActiveDocument

# 1. Delete the "Meta description" paragraph (the bold "Meta description"
#    label followed by the intro-text run for the review).
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Meta description*") {
        $p.Range.Delete()
        break
    }
}

# 2. Insert a new bold paragraph "Play FashionTV Highlife for Free - Review 2021"
#    right before the final (italic) paragraph at the end of the document.
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$insertionPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play FashionTV Highlife for Free - Review 2021</w:t></w:r></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertionPoint.InsertXML($xml)

# InsertXML leaves behind an extra empty paragraph mark used only to force
# the paragraph break; remove it (it is the paragraph consisting solely of
# its end-of-paragraph mark, i.e. Range length == 1).
$count2 = $d.Paragraphs.Count
for ($i = $count2; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if (($p.Range.End - $p.Range.Start) -eq 1) {
        $p.Range.Delete()
        break
    }
}

# 3. Replace the text of the final (italic) paragraph with the new copy.
$old = "Create a feature image for FashionTV Highlife slot game that captures the luxurious and glamorous atmosphere of the game. The image should be in a cartoon style, depicting a happy Maya warrior with glasses, surrounded by the symbols of the game, such as a car, yacht, gold watch, and the three fantastic models. The background should feature the FashionTV Highlife logo and a vibrant city at night with bright lights shining. Make sure to include the Special Nudging HP1 and Free Spins function icons in the image, emphasizing the special features of the game. Overall, the feature image should convey the excitement of winning high sums while enjoying the extravagance and luxury of this online slot game."
$new = "Read our impartial review of FashionTV Highlife online slot, learn how to play and where to play for free in 2021."
$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)

Write-Output "edit applied"
